$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 22.700661
$ws.Cells.Item(2, 8).Value = 68.10198299999999
$ws.Cells.Item(2, 9).Value = 0.08615268874617349
$ws.Cells.Item(2, 10).Value = 0.08615268874617349
$ws.Cells.Item(2, 13).Value = 255.0443116666667
$ws.Cells.Item(2, 14).Value = 765.132935
$ws.Cells.Item(2, 15).Value = 0.863617428561108
$ws.Cells.Item(2, 16).Value = 0.8636174285611079
$ws.Cells.Item(2, 17).Value = 5789.674459123344
$ws.Cells.Item(2, 18).Value = 52107.07013211009
$ws.Cells.Item(2, 19).Value = 0.07440296351859586
$ws.Cells.Item(2, 20).Value = 0.07440296351859585
$ws.Cells.Item(3, 7).Value = 22.700661
$ws.Cells.Item(3, 8).Value = 68.10198299999999
$ws.Cells.Item(3, 9).Value = 0.08615268874617349
$ws.Cells.Item(3, 10).Value = 0.08615268874617349
$ws.Cells.Item(3, 13).Value = 0.8952453333333334
$ws.Cells.Item(3, 15).Value = 0.003031431940796009
$ws.Cells.Item(3, 16).Value = 0.003031431940796009
$ws.Cells.Item(3, 17).Value = 20.322660823832
$ws.Cells.Item(3, 18).Value = 182.903947414488
$ws.Cells.Item(3, 19).Value = 0.0002611660124506071
$ws.Cells.Item(3, 20).Value = 0.0002611660124506071
$ws.Cells.Item(4, 7).Value = 22.700661
$ws.Cells.Item(4, 8).Value = 68.10198299999999
$ws.Cells.Item(4, 9).Value = 0.08615268874617349
$ws.Cells.Item(4, 10).Value = 0.08615268874617349
$ws.Cells.Item(4, 13).Value = 7.050555333333333
$ws.Cells.Item(4, 14).Value = 21.151666
$ws.Cells.Item(4, 15).Value = 0.02387421396349043
$ws.Cells.Item(4, 16).Value = 0.02387421396349043
$ws.Cells.Item(4, 17).Value = 160.052266483742
$ws.Cells.Item(4, 18).Value = 1440.470398353678
$ws.Cells.Item(4, 19).Value = 0.00205682772465614
$ws.Cells.Item(4, 20).Value = 0.00205682772465614
$ws.Cells.Item(5, 7).Value = 22.700661
$ws.Cells.Item(5, 8).Value = 68.10198299999999
$ws.Cells.Item(5, 9).Value = 0.08615268874617349
$ws.Cells.Item(5, 10).Value = 0.08615268874617349
$ws.Cells.Item(5, 13).Value = 32.33082866666666
$ws.Cells.Item(5, 14).Value = 96.99248599999999
$ws.Cells.Item(5, 15).Value = 0.1094769255346056
$ws.Cells.Item(5, 16).Value = 0.1094769255346056
$ws.Cells.Item(5, 17).Value = 733.9311814110818
$ws.Cells.Item(5, 18).Value = 6605.380632699736
$ws.Cells.Item(5, 19).Value = 0.009431731490470892
$ws.Cells.Item(5, 20).Value = 0.009431731490470892
$ws.Cells.Item(6, 9).Value = 0.5030288587986086
$ws.Cells.Item(6, 10).Value = 0.5030288587986087
$ws.Cells.Item(6, 13).Value = 255.0443116666667
$ws.Cells.Item(6, 14).Value = 765.132935
$ws.Cells.Item(6, 15).Value = 0.863617428561108
$ws.Cells.Item(6, 16).Value = 0.8636174285611079
$ws.Cells.Item(6, 17).Value = 33804.78750429738
$ws.Cells.Item(6, 18).Value = 304243.0875386765
$ws.Cells.Item(6, 19).Value = 0.434424489527683
$ws.Cells.Item(6, 20).Value = 0.4344244895276831
$ws.Cells.Item(7, 9).Value = 0.5030288587986086
$ws.Cells.Item(7, 10).Value = 0.5030288587986087
$ws.Cells.Item(7, 13).Value = 0.8952453333333334
$ws.Cells.Item(7, 15).Value = 0.003031431940796009
$ws.Cells.Item(7, 16).Value = 0.003031431940796009
$ws.Cells.Item(7, 19).Value = 0.001524897749704268
$ws.Cells.Item(7, 20).Value = 0.001524897749704268
$ws.Cells.Item(8, 9).Value = 0.5030288587986086
$ws.Cells.Item(8, 10).Value = 0.5030288587986087
$ws.Cells.Item(8, 13).Value = 7.050555333333333
$ws.Cells.Item(8, 14).Value = 21.151666
$ws.Cells.Item(8, 15).Value = 0.02387421396349043
$ws.Cells.Item(8, 16).Value = 0.02387421396349043
$ws.Cells.Item(8, 17).Value = 934.514176274312
$ws.Cells.Item(8, 18).Value = 8410.627586468809
$ws.Cells.Item(8, 19).Value = 0.0120094186047684
$ws.Cells.Item(8, 20).Value = 0.0120094186047684
$ws.Cells.Item(9, 9).Value = 0.5030288587986086
$ws.Cells.Item(9, 10).Value = 0.5030288587986087
$ws.Cells.Item(9, 13).Value = 32.33082866666666
$ws.Cells.Item(9, 14).Value = 96.99248599999999
$ws.Cells.Item(9, 15).Value = 0.1094769255346056
$ws.Cells.Item(9, 16).Value = 0.1094769255346056
$ws.Cells.Item(9, 17).Value = 4285.281980109167
$ws.Cells.Item(9, 18).Value = 38567.5378209825
$ws.Cells.Item(9, 19).Value = 0.05507005291645293
$ws.Cells.Item(9, 20).Value = 0.05507005291645294
$ws.Cells.Item(10, 7).Value = 41.94534433333333
$ws.Cells.Item(10, 8).Value = 125.836033
$ws.Cells.Item(10, 9).Value = 0.159189381961201
$ws.Cells.Item(10, 10).Value = 0.159189381961201
$ws.Cells.Item(10, 13).Value = 255.0443116666667
$ws.Cells.Item(10, 14).Value = 765.132935
$ws.Cells.Item(10, 15).Value = 0.863617428561108
$ws.Cells.Item(10, 16).Value = 0.8636174285611079
$ws.Cells.Item(10, 17).Value = 10697.92147311632
$ws.Cells.Item(10, 18).Value = 96281.29325804685
$ws.Cells.Item(10, 19).Value = 0.1374787247035645
$ws.Cells.Item(10, 20).Value = 0.1374787247035644
$ws.Cells.Item(11, 7).Value = 41.94534433333333
$ws.Cells.Item(11, 8).Value = 125.836033
$ws.Cells.Item(11, 9).Value = 0.159189381961201
$ws.Cells.Item(11, 10).Value = 0.159189381961201
$ws.Cells.Item(11, 13).Value = 0.8952453333333334
$ws.Cells.Item(11, 15).Value = 0.003031431940796009
$ws.Cells.Item(11, 16).Value = 0.003031431940796009
$ws.Cells.Item(11, 17).Value = 37.55137376947645
$ws.Cells.Item(11, 18).Value = 337.962363925288
$ws.Cells.Item(11, 19).Value = 0.0004825717771127606
$ws.Cells.Item(11, 20).Value = 0.0004825717771127606
$ws.Cells.Item(12, 7).Value = 41.94534433333333
$ws.Cells.Item(12, 8).Value = 125.836033
$ws.Cells.Item(12, 9).Value = 0.159189381961201
$ws.Cells.Item(12, 10).Value = 0.159189381961201
$ws.Cells.Item(12, 13).Value = 7.050555333333333
$ws.Cells.Item(12, 14).Value = 21.151666
$ws.Cells.Item(12, 15).Value = 0.02387421396349043
$ws.Cells.Item(12, 16).Value = 0.02387421396349043
$ws.Cells.Item(12, 17).Value = 295.7379711978864
$ws.Cells.Item(12, 18).Value = 2661.641740780978
$ws.Cells.Item(12, 19).Value = 0.003800521365657517
$ws.Cells.Item(12, 20).Value = 0.003800521365657516
$ws.Cells.Item(13, 7).Value = 41.94534433333333
$ws.Cells.Item(13, 8).Value = 125.836033
$ws.Cells.Item(13, 9).Value = 0.159189381961201
$ws.Cells.Item(13, 10).Value = 0.159189381961201
$ws.Cells.Item(13, 13).Value = 32.33082866666666
$ws.Cells.Item(13, 14).Value = 96.99248599999999
$ws.Cells.Item(13, 15).Value = 0.1094769255346056
$ws.Cells.Item(13, 16).Value = 0.1094769255346056
$ws.Cells.Item(13, 17).Value = 1356.127741005337
$ws.Cells.Item(13, 18).Value = 12205.14966904803
$ws.Cells.Item(13, 19).Value = 0.0174275641148663
$ws.Cells.Item(13, 20).Value = 0.01742756411486629
$ws.Cells.Item(14, 7).Value = 66.302588
$ws.Cells.Item(14, 8).Value = 198.907764
$ws.Cells.Item(14, 9).Value = 0.2516290704940168
$ws.Cells.Item(14, 10).Value = 0.2516290704940168
$ws.Cells.Item(14, 13).Value = 255.0443116666667
$ws.Cells.Item(14, 14).Value = 765.132935
$ws.Cells.Item(14, 15).Value = 0.863617428561108
$ws.Cells.Item(14, 16).Value = 0.8636174285611079
$ws.Cells.Item(14, 17).Value = 16910.09791817859
$ws.Cells.Item(14, 18).Value = 152190.8812636073
$ws.Cells.Item(14, 19).Value = 0.2173112508112646
$ws.Cells.Item(14, 20).Value = 0.2173112508112646
$ws.Cells.Item(15, 7).Value = 66.302588
$ws.Cells.Item(15, 8).Value = 198.907764
$ws.Cells.Item(15, 9).Value = 0.2516290704940168
$ws.Cells.Item(15, 10).Value = 0.2516290704940168
$ws.Cells.Item(15, 13).Value = 0.8952453333333334
$ws.Cells.Item(15, 15).Value = 0.003031431940796009
$ws.Cells.Item(15, 16).Value = 0.003031431940796009
$ws.Cells.Item(15, 17).Value = 59.35708249492267
$ws.Cells.Item(15, 18).Value = 534.213742454304
$ws.Cells.Item(15, 19).Value = 0.0007627964015283731
$ws.Cells.Item(15, 20).Value = 0.0007627964015283731
$ws.Cells.Item(16, 7).Value = 66.302588
$ws.Cells.Item(16, 8).Value = 198.907764
$ws.Cells.Item(16, 9).Value = 0.2516290704940168
$ws.Cells.Item(16, 10).Value = 0.2516290704940168
$ws.Cells.Item(16, 13).Value = 7.050555333333333
$ws.Cells.Item(16, 14).Value = 21.151666
$ws.Cells.Item(16, 15).Value = 0.02387421396349043
$ws.Cells.Item(16, 16).Value = 0.02387421396349043
$ws.Cells.Item(16, 17).Value = 467.4700654372026
$ws.Cells.Item(16, 18).Value = 4207.230588934824
$ws.Cells.Item(16, 19).Value = 0.006007446268408375
$ws.Cells.Item(16, 20).Value = 0.006007446268408374
$ws.Cells.Item(17, 7).Value = 66.302588
$ws.Cells.Item(17, 8).Value = 198.907764
$ws.Cells.Item(17, 9).Value = 0.2516290704940168
$ws.Cells.Item(17, 10).Value = 0.2516290704940168
$ws.Cells.Item(17, 13).Value = 32.33082866666666
$ws.Cells.Item(17, 14).Value = 96.99248599999999
$ws.Cells.Item(17, 15).Value = 0.1094769255346056
$ws.Cells.Item(17, 16).Value = 0.1094769255346056
$ws.Cells.Item(17, 17).Value = 2143.617612784589
$ws.Cells.Item(17, 18).Value = 19292.5585150613
$ws.Cells.Item(17, 19).Value = 0.02754757701281552
$ws.Cells.Item(17, 20).Value = 0.02754757701281551
